$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 185, shifting rows 185:234 down to 186:235
$ws.Rows("185:185").Insert()

# Populate the newly inserted row 185 with the new data point
$ws.Range("A185").Value = 9
$ws.Range("B185").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C185").Value = "Metropolitana"
$ws.Range("D185").Value = 44588
$ws.Range("E185").Value = 13
$ws.Range("F185").Value = 100112001
$ws.Range("G185").Value = "Berenjena"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 70
$ws.Range("K185").Value = 13000
$ws.Range("L185").Value = 14000
$ws.Range("M185").Value = 13500
$ws.Range("N185").Value = "`$/caja 60 unidades"
$ws.Range("O185").Value = "Región Metropolitana"
$ws.Range("P185").Value = 225
$ws.Range("Q185").Value = 60
$ws.Range("R185").Value = "Hortaliza"
